# Daily attendance processing - 2025-12-12 23:49:53
# Normalizes the "Recorded By" (column G) text so that "System" is always
# listed before the human recorder's e-mail address, and so that the
# duplicate lowercase "system" token in the backup account rows comes
# after "System" instead of before it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "backup@backdoor.com, system, System") {
        $cell.Value2 = "backup@backdoor.com, System, system"
    }
}
